$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1004.13043
$ws.Range("I107").Value = 793.4211
$ws.Range("K107").Value = 793.4211
$ws.Range("M107").Value = 1126.5789
$ws.Range("H111").Value = 1306.091
$ws.Range("I111").Value = 468.44446
$ws.Range("J111").Value = 5075.5
$ws.Range("K111").Value = 1405.33338
$ws.Range("L111").Value = 15226.5
$ws.Range("M111").Value = 1661.66662
$ws.Range("N111").Value = -21360.5
$ws.Range("H112").Value = 727894.7
$ws.Range("I112").Value = 2197.5
$ws.Range("J112").Value = 839540.4
$ws.Range("K112").Value = 6592.5
$ws.Range("L112").Value = 2518621.2
$ws.Range("M112").Value = -5484.5
$ws.Range("N112").Value = -2520837.2
$ws.Range("H133").Value = 79500
$ws.Range("J133").Value = 79500
$ws.Range("L133").Value = 79500
$ws.Range("N133").Value = -89620
$ws.Range("H137").Value = 1605.3429
$ws.Range("I137").Value = 836.15
$ws.Range("K137").Value = 2508.45
$ws.Range("M137").Value = 41.55000000000018
$ws.Range("H138").Value = 2901.5393
$ws.Range("J138").Value = 3320.014
$ws.Range("L138").Value = 9960.042000000001
$ws.Range("N138").Value = -20240.042

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3005711.2
$ws.Range("I32").Value = 3130443.5
$ws.Range("K32").Value = 3130443.5
$ws.Range("M32").Value = -3130156.5
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("H45").Value = 3522.0667
$ws.Range("I45").Value = 2890.889
$ws.Range("K45").Value = 2890.889
$ws.Range("M45").Value = -2513.889
$ws.Range("H61").Value = 2719.111
$ws.Range("I61").Value = 2559.125
$ws.Range("J61").Value = 3999
$ws.Range("K61").Value = 2559.125
$ws.Range("L61").Value = 3999
$ws.Range("M61").Value = -2347.125
$ws.Range("N61").Value = -4423
$ws.Range("H74").Value = 1940.5714
$ws.Range("I74").Value = 1583.7059
$ws.Range("J74").Value = 2749.4666
$ws.Range("K74").Value = 1583.7059
$ws.Range("L74").Value = 2749.4666
$ws.Range("M74").Value = -709.7058999999999
$ws.Range("N74").Value = -4497.4666
$ws.Range("H77").Value = 1940.5714
$ws.Range("I77").Value = 1583.7059
$ws.Range("J77").Value = 2749.4666
$ws.Range("K77").Value = 7918.5295
$ws.Range("L77").Value = 13747.333
$ws.Range("M77").Value = -3550.5295
$ws.Range("N77").Value = -22483.333
$ws.Range("H122").Value = 1568.9048
$ws.Range("I122").Value = 1126.1177
$ws.Range("K122").Value = 3378.3531
$ws.Range("M122").Value = -928.3531000000003
$ws.Range("H132").Value = 22729588
$ws.Range("I132").Value = 2279.353
$ws.Range("J132").Value = 100002440
$ws.Range("K132").Value = 6838.059
$ws.Range("L132").Value = 300007320
$ws.Range("M132").Value = -4308.059
$ws.Range("N132").Value = -300012380
$ws.Range("H136").Value = 2719.111
$ws.Range("I136").Value = 2559.125
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 7677.375
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -5127.375
$ws.Range("N136").Value = -17097
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 69107.3
$ws.Range("I105").Value = 2549.55
$ws.Range("K105").Value = 2549.55
$ws.Range("M105").Value = -802.5500000000002
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H132").Value = 99999
$ws.Range("J132").Value = 99999
$ws.Range("L132").Value = 99999
$ws.Range("N132").Value = -110119
$ws.Range("H134").Value = 21607808
$ws.Range("J134").Value = 166669170
$ws.Range("L134").Value = 500007510
$ws.Range("N134").Value = -500012580

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1741.4546
$ws.Range("I16").Value = 1458.25
$ws.Range("J16").Value = 2496.6667
$ws.Range("K16").Value = 1458.25
$ws.Range("L16").Value = 2496.6667
$ws.Range("M16").Value = -1171.25
$ws.Range("N16").Value = -3070.6667
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = 150
$ws.Range("H32").Value = 11626.25
$ws.Range("I32").Value = 11626.25
$ws.Range("K32").Value = 11626.25
$ws.Range("M32").Value = -11310.25
$ws.Range("H99").Value = 2755
$ws.Range("I99").Value = 1981.8572
$ws.Range("J99").Value = 3528.1428
$ws.Range("K99").Value = 1981.8572
$ws.Range("L99").Value = 3528.1428
$ws.Range("M99").Value = -483.8571999999999
$ws.Range("N99").Value = -6524.1428
$ws.Range("H105").Value = 2277.8667
$ws.Range("I105").Value = 1537.7273
$ws.Range("K105").Value = 1537.7273
$ws.Range("M105").Value = 209.2727
$ws.Range("H113").Value = 1741.4546
$ws.Range("I113").Value = 1458.25
$ws.Range("J113").Value = 2496.6667
$ws.Range("K113").Value = 1458.25
$ws.Range("L113").Value = 2496.6667
$ws.Range("M113").Value = 711.75
$ws.Range("N113").Value = -6836.6667
$ws.Range("H126").Value = 2755
$ws.Range("I126").Value = 1981.8572
$ws.Range("J126").Value = 3528.1428
$ws.Range("K126").Value = 5945.571599999999
$ws.Range("L126").Value = 10584.4284
$ws.Range("M126").Value = -3475.571599999999
$ws.Range("N126").Value = -15524.4284
$ws.Range("H134").Value = 3572667.8
$ws.Range("I134").Value = 1277.6086
$ws.Range("K134").Value = 3832.8258
$ws.Range("M134").Value = -1297.8258

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3366.923
$ws.Range("I129").Value = 993
$ws.Range("J129").Value = 3564.75
$ws.Range("K129").Value = 2979
$ws.Range("L129").Value = 10694.25
$ws.Range("M129").Value = 2021
$ws.Range("N129").Value = -20694.25
$ws.Range("H137").Value = 7281.2144
$ws.Range("I137").Value = 2701
$ws.Range("J137").Value = 11861.429
$ws.Range("K137").Value = 8103
$ws.Range("L137").Value = 35584.287
$ws.Range("M137").Value = -3003
$ws.Range("N137").Value = -45784.287

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 201.72223
$ws.Range("I2").Value = 74.42856999999999
$ws.Range("J2").Value = 647.25
$ws.Range("K2").Value = 74.42856999999999
$ws.Range("L2").Value = 647.25
$ws.Range("M2").Value = 38.57143000000001
$ws.Range("N2").Value = -873.25
$ws.Range("H102").Value = 6013.5
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 6013.5
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 6013.5
$ws.Range("N102").Value = -9257.5
$ws.Range("H122").Value = 6317.3335
$ws.Range("I122").Value = 25000
$ws.Range("J122").Value = 3982
$ws.Range("K122").Value = 75000
$ws.Range("L122").Value = 11946
$ws.Range("M122").Value = -72550
$ws.Range("N122").Value = -16846
$ws.Range("H132").Value = 1361.3125
$ws.Range("J132").Value = 1999.25
$ws.Range("L132").Value = 5997.75
$ws.Range("N132").Value = -11057.75
$ws.Range("H133").Value = 90337.5
$ws.Range("J133").Value = 90337.5
$ws.Range("L133").Value = 90337.5
$ws.Range("N133").Value = -100457.5
$ws.Range("H136").Value = 76397.2
$ws.Range("J136").Value = 76397.2
$ws.Range("L136").Value = 229191.6
$ws.Range("N136").Value = -234291.6
$ws.Range("M102").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3682.8572
$ws.Range("I7").Value = 3470.6
$ws.Range("J7").Value = 4213.5
$ws.Range("K7").Value = 3470.6
$ws.Range("L7").Value = 4213.5
$ws.Range("M7").Value = -3358.6
$ws.Range("N7").Value = -4437.5
$ws.Range("H40").Value = 7056.75
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 7056.75
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 7056.75
$ws.Range("N40").Value = -7328.75
$ws.Range("H53").Value = 29999
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 29999
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 29999
$ws.Range("N53").Value = -31035
$ws.Range("H126").Value = 3682.8572
$ws.Range("I126").Value = 3470.6
$ws.Range("J126").Value = 4213.5
$ws.Range("K126").Value = 10411.8
$ws.Range("L126").Value = 12640.5
$ws.Range("M126").Value = -7941.799999999999
$ws.Range("N126").Value = -17580.5
$ws.Range("H132").Value = 2236.1667
$ws.Range("I132").Value = 2236.1667
$ws.Range("K132").Value = 6708.500100000001
$ws.Range("M132").Value = -4178.500100000001
$ws.Range("H136").Value = 28576560
$ws.Range("I136").Value = 4611.483
$ws.Range("J136").Value = 166674300
$ws.Range("K136").Value = 13834.449
$ws.Range("L136").Value = 500022900
$ws.Range("M136").Value = -11284.449
$ws.Range("N136").Value = -500028000
$ws.Range("M40").ClearContents()
$ws.Range("M53").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 505
$ws.Range("I29").Value = 505
$ws.Range("K29").Value = 505
$ws.Range("M29").Value = -215
$ws.Range("H38").Value = 9668.5
$ws.Range("J38").Value = 9822.200000000001
$ws.Range("L38").Value = 9822.200000000001
$ws.Range("N38").Value = -10768.2
$ws.Range("H100").Value = 1407.8572
$ws.Range("I100").Value = 963.25
$ws.Range("J100").Value = 2000.6666
$ws.Range("K100").Value = 1926.5
$ws.Range("L100").Value = 4001.3332
$ws.Range("M100").Value = -1385.5
$ws.Range("N100").Value = -5083.3332
$ws.Range("H126").Value = 1759.6111
$ws.Range("I126").Value = 1392.2354
$ws.Range("J126").Value = 8005
$ws.Range("K126").Value = 4176.706200000001
$ws.Range("L126").Value = 24015
$ws.Range("M126").Value = -1706.706200000001
$ws.Range("N126").Value = -28955
$ws.Range("H136").Value = 1039.0646
$ws.Range("I136").Value = 808.38464
$ws.Range("J136").Value = 2238.6
$ws.Range("K136").Value = 2425.15392
$ws.Range("L136").Value = 6715.799999999999
$ws.Range("M136").Value = 124.8460800000003
$ws.Range("N136").Value = -11815.8
